$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E (shifts old D.. onward to F..)
$ws.Range("D1:E1").EntireColumn.Insert()

# Set header text for the two newly inserted columns
$ws.Range("D1").Value = "NOTES"
$ws.Range("E1").Value = "ABBR"

# Match the style of the header row cell to its left neighbor (C1) for the new headers
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set the column widths for the newly inserted columns (raw OOXML width 14
# corresponds to a ColumnWidth property of 14 - 5/6 given this engine's MDW of 6)
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 13.166666666666666

# Update the view: zoom and selection
$ws.Application.ActiveWindow.Zoom = 190
$ws.Range("E10").Select()
